$wb = $excel.ActiveWorkbook

# Remember which sheet is active so we can restore it - adding a sheet
# otherwise leaves the newly created sheet selected/active.
$prevActiveName = $wb.ActiveSheet.Name

# Add a brand-new worksheet immediately after the last existing sheet, so it
# lands at the very end of the workbook (tab order), and name it "qwe".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "qwe"

# Header row content (single row of column headings, A1:F1).
$ws.Range("A1").Value = "Outdoor Model"
$ws.Range("B1").Value = "Outdoor Quantity"
$ws.Range("C1").Value = "Outdoor Serial(s)"
$ws.Range("D1").Value = "Indoor Model"
$ws.Range("E1").Value = "Indoor Quantity"
$ws.Range("F1").Value = "Indoor Serial(s)"

# Formatting to match the other sheets' header style: bold font, thin box
# border around each cell, centered horizontally and aligned to top.
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Match the page margins used throughout the rest of the workbook.
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(1)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$ws.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$ws.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

$ws.Range("A1").Select() | Out-Null

# Restore the previously active sheet/tab so the new sheet is appended
# without disturbing which tab is shown as selected.
$wb.Worksheets.Item($prevActiveName).Activate() | Out-Null
